# 2017-4-13 log4net 与 Quartz实例签入
# Adds two new task rows (Quartz info filled into existing row 13, and a
# brand-new row 16 for Log4net), matching date/status cell formatting from
# existing rows, and moves the active selection to F20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: fill in Quartz task's dates / status / project columns ------
# Row 13 already has A13=12 ("序号") and B13 = "Quartz定时任务框架实例".
# Copy date-formatting (style s=1) from an existing date cell, and the
# "已完成" status formatting (style s=2, green fill) from an existing one,
# then stamp in the values.
$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(13, 3).Value = 42835          # 2017-04-10

$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(13, 4).Value = 42838          # 2017-04-13

$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(13, 5).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(13, 5).Value = "已完成"

$ws.Cells.Item(13, 6).Value = "QuartzTest"

# --- Row 16 (new): Log4net task -------------------------------------------
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "Log4net记录日志实例"

$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(16, 3).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(16, 3).Value = 42838          # 2017-04-13

$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(16, 4).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(16, 4).Value = 42838          # 2017-04-13

$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(16, 5).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(16, 5).Value = "已完成"

$ws.Cells.Item(16, 6).Value = "Log4netTest"

# --- Move the selection, as recorded in the saved view -------------------
$ws.Range("F20").Select() | Out-Null

Write-Host "edit complete"
